# Commit: "added kiana topic 1, section 1"
#
# The existing data block for Kiana / "Section 1, Expressions" only has
# Topic numbers 2-9 (Topic 1 is missing). This inserts a new row right
# before the current Kiana block (worksheet row 39) holding Kiana's
# missing Topic 1 entry, pushing every following row down by one.
#
# The new row's Time/Time2 cells reuse the same placeholder label that
# was already sitting in (old) row 39 ("33:16"), matching exactly what
# happened in the authored edit (a row insert/copy-down, not new data).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the first Kiana / Section 1 row, shifting
# rows 39:130 down to 40:131.
$ws.Rows(39).Insert()

# Populate the freshly inserted row 39 with Kiana's Topic 1 entry.
$ws.Range("A39").Value = "Kiana"
$ws.Range("B39").Value = "Section 1, Expressions"
$ws.Range("C39").Value = 1

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "33:16"

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "33:16"

# Restore the selection / scroll position recorded at save time.
$ws.Range("F39").Select()
$excel.ActiveWindow.ScrollRow = 32
$excel.ActiveWindow.ScrollColumn = 1
